$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.789499
$ws.Range("H2").Value = 5.368497
$ws.Range("I2").Value = 0.01244533957901722
$ws.Range("J2").Value = 0.01244533957901722
$ws.Range("M2").Value = 1.975599
$ws.Range("N2").Value = 5.926797000000001
$ws.Range("O2").Value = 0.2862071854613167
$ws.Range("P2").Value = 0.2862071854613167
$ws.Range("Q2").Value = 3.535332434901
$ws.Range("R2").Value = 31.817991914109
$ws.Range("S2").Value = 0.003561945613020847
$ws.Range("T2").Value = 0.003561945613020846
$ws.Range("G3").Value = 1.789499
$ws.Range("H3").Value = 5.368497
$ws.Range("I3").Value = 0.01244533957901722
$ws.Range("J3").Value = 0.01244533957901722
$ws.Range("O3").Value = 0.6582911054219697
$ws.Range("P3").Value = 0.6582911054219697
$ws.Range("Q3").Value = 8.131444683521664
$ws.Range("R3").Value = 73.18300215169499
$ws.Range("S3").Value = 0.008192656348823036
$ws.Range("T3").Value = 0.008192656348823036
$ws.Range("G4").Value = 1.789499
$ws.Range("H4").Value = 5.368497
$ws.Range("I4").Value = 0.01244533957901722
$ws.Range("J4").Value = 0.01244533957901722
$ws.Range("O4").Value = 0.05550170911671371
$ws.Range("P4").Value = 0.05550170911671371
$ws.Range("Q4").Value = 0.685576751389
$ws.Range("R4").Value = 6.170190762500999
$ws.Range("S4").Value = 0.000690737617173338
$ws.Range("T4").Value = 0.000690737617173338
$ws.Range("H5").Value = 311.722962
$ws.Range("I5").Value = 0.7226413867171911
$ws.Range("J5").Value = 0.7226413867171912
$ws.Range("M5").Value = 1.975599
$ws.Range("N5").Value = 5.926797000000001
$ws.Range("O5").Value = 0.2862071854613167
$ws.Range("P5").Value = 0.2862071854613167
$ws.Range("Q5").Value = 205.279857334746
$ws.Range("R5").Value = 1847.518716012714
$ws.Range("S5").Value = 0.2068251573901902
$ws.Range("T5").Value = 0.2068251573901902
$ws.Range("H6").Value = 311.722962
$ws.Range("I6").Value = 0.7226413867171911
$ws.Range("J6").Value = 0.7226413867171912
$ws.Range("O6").Value = 0.6582911054219697
$ws.Range("P6").Value = 0.6582911054219697
$ws.Range("Q6").Value = 472.1541284434966
$ws.Range("R6").Value = 4249.38715599147
$ws.Range("S6").Value = 0.4757083972857248
$ws.Range("T6").Value = 0.4757083972857249
$ws.Range("H7").Value = 311.722962
$ws.Range("I7").Value = 0.7226413867171911
$ws.Range("J7").Value = 0.7226413867171912
$ws.Range("O7").Value = 0.05550170911671371
$ws.Range("P7").Value = 0.05550170911671371
$ws.Range("Q7").Value = 39.808165231594
$ws.Range("R7").Value = 358.273487084346
$ws.Range("S7").Value = 0.04010783204127617
$ws.Range("T7").Value = 0.04010783204127617
$ws.Range("I8").Value = 0.2649132737037916
$ws.Range("J8").Value = 0.2649132737037916
$ws.Range("M8").Value = 1.975599
$ws.Range("N8").Value = 5.926797000000001
$ws.Range("O8").Value = 0.2862071854613167
$ws.Range("P8").Value = 0.2862071854613167
$ws.Range("Q8").Value = 75.253590552069
$ws.Range("R8").Value = 677.282314968621
$ws.Range("S8").Value = 0.07582008245810563
$ws.Range("T8").Value = 0.07582008245810562
$ws.Range("I9").Value = 0.2649132737037916
$ws.Range("J9").Value = 0.2649132737037916
$ws.Range("O9").Value = 0.6582911054219697
$ws.Range("P9").Value = 0.6582911054219697
$ws.Range("S9").Value = 0.1743900517874218
$ws.Range("T9").Value = 0.1743900517874218
$ws.Range("I10").Value = 0.2649132737037916
$ws.Range("J10").Value = 0.2649132737037916
$ws.Range("O10").Value = 0.05550170911671371
$ws.Range("P10").Value = 0.05550170911671371
$ws.Range("S10").Value = 0.0147031394582642
$ws.Range("T10").Value = 0.0147031394582642
